{"js": "// Remove the trailing site-footer paragraphs from the course page:\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n// plus the now-redundant blank paragraph that used to separate the\n// copyright line from the trailing page-break paragraph. The blank\n// paragraph that immediately follows \"LOB1012: Estat\u00edstica (Requisito\n// fraco)\" is left untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer text paragraphs by their content so the script\n// is resilient to exact indices.\nconst jupiterIdx = items.findIndex(\n  (p) => p.text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n);\nif (jupiterIdx === -1) {\n  throw new Error('Paragraph \"Ver no Jupiter Salvar em pdf Salvar em docx\" not found.');\n}\n\nconst copyrightIdx = jupiterIdx + 1;\nif (items[copyrightIdx].text.indexOf(\"Contact: luizeleno@usp.br\") === -1) {\n  throw new Error(\"Expected copyright paragraph right after the Jupiter paragraph.\");\n}\n\n// The blank paragraph right after the copyright line is removed along\n// with the two text paragraphs (it duplicates the blank paragraph that\n// already sits right after the \"LOB1012\" line).\nconst trailingBlankIdx = copyrightIdx + 1;\nif (items[trailingBlankIdx].text !== \"\") {\n  throw new Error(\"Expected a blank paragraph after the copyright paragraph.\");\n}\n\n// Delete from the highest index down so earlier indices stay valid.\nitems[trailingBlankIdx].delete();\nitems[copyrightIdx].delete();\nitems[jupiterIdx].delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing site-footer paragraphs from the course page:\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n# plus the now-redundant blank paragraph that used to separate the\n# copyright line from the trailing page-break paragraph. The blank\n# paragraph that immediately follows \"LOB1012: Estatistica (Requisito\n# fraco)\" is left untouched.\n\n$d = $word.ActiveDocument\n\n$jupiterIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -like \"Ver no Jupiter*\") {\n        $jupiterIdx = $i\n        break\n    }\n}\nif ($jupiterIdx -eq -1) {\n    throw \"Could not find the 'Ver no Jupiter...' paragraph.\"\n}\n\n$copyrightIdx = $jupiterIdx + 1\n$copyrightText = $d.Paragraphs.Item($copyrightIdx).Range.Text\nif ($copyrightText -notlike \"*Contact: luizeleno@usp.br*\") {\n    throw \"Expected copyright paragraph right after the Jupiter paragraph.\"\n}\n\n$trailingBlankIdx = $copyrightIdx + 1\n$blankText = $d.Paragraphs.Item($trailingBlankIdx).Range.Text.Trim()\nif ($blankText.Length -ne 0) {\n    throw \"Expected a blank paragraph after the copyright paragraph.\"\n}\n\n# Delete from the highest index down so earlier indices stay valid.\n$d.Paragraphs.Item($trailingBlankIdx).Range.Delete()\n$d.Paragraphs.Item($copyrightIdx).Range.Delete()\n$d.Paragraphs.Item($jupiterIdx).Range.Delete()\n"}
